# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.77 = 51576.66 pesos`n✅ 51576.66 pesos = 12.68 = 979.03 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 78.33
$wsTasas.Range("O10").Value = 4040
$wsTasas.Range("N12").Value = 4067
$wsTasas.Range("O12").Value = 77.2
